# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet right after "总计" (the existing
# "2022-Q2" / "2022-Q1" / "2021-Q4" sheets simply shift one tab to the
# right, unchanged), fills it with the new quarter's fund-holding table,
# and records the new quarter as the top data row of the "总计" summary
# sheet (pushing its previous rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Reuse the header row + column-A formatting from the existing "2022-Q2"
# sheet (same table layout) via Copy, which carries the cell style along
# with the value - more reliable than the `.Style` property here.
$q2Sheet.Range("A1:H1").Copy($q3.Range("A1:H1"))
$q2Sheet.Range("A2:A2").Copy($q3.Range("A2:A12"))

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Text-typed columns (B, D, E, F, G) must stay text - many values are
# numeric-looking strings ("23.45") or codes with leading zeros ("014085")
# that Excel would otherwise silently coerce to numbers, corrupting them.
$q3.Range("B2:B12").NumberFormat = "@"
$q3.Range("D2:G12").NumberFormat = "@"

# A column (row index, 0-based) and H column (rank) are genuine numbers.
$q3Data = @(
    @("519692", "交银成长混合A", "23.45", "76.71", "2.08", "0.4878", "10"),
    @("519694", "交银蓝筹混合", "16.10", "78.09", "2.07", "0.3333", "10"),
    @("506007", "广发科创板两年定开混合", "5.01", "94.25", "6.56", "0.3287", "3"),
    @("110002", "易方达策略成长混合", "11.73", "90.54", "2.77", "0.3249", "8"),
    @("112002", "易方达策略成长二号混合", "9.94", "91.00", "2.76", "0.2743", "8"),
    @("166801", "浙商聚潮新思维混合A", "1.89", "78.76", "3.51", "0.0663", "9"),
    @("014085", "浙商聚潮新思维混合C", "0.85", "78.76", "3.51", "0.0298", "9"),
    @("015373", "浙商智选新兴产业混合A", "0.70", "92.12", "4.16", "0.0291", "7"),
    @("015374", "浙商智选新兴产业混合C", "0.27", "92.12", "4.16", "0.0112", "7"),
    @("006538", "东海核心价值精选混合", "0.19", "60.22", "2.33", "0.0044", "7"),
    @("960016", "交银成长混合H", "0.16", "76.71", "2.08", "0.0033", "10")
)

$row = 2
foreach ($r in $q3Data) {
    $q3.Cells.Item($row, 1).Value = $row - 2
    $q3.Cells.Item($row, 2).Value = $r[0]
    $q3.Cells.Item($row, 3).Value = $r[1]
    $q3.Cells.Item($row, 4).Value = $r[2]
    $q3.Cells.Item($row, 5).Value = $r[3]
    $q3.Cells.Item($row, 6).Value = $r[4]
    $q3.Cells.Item($row, 7).Value = $r[5]
    $q3.Cells.Item($row, 8).Value = [int]$r[6]
    $row++
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new top data row for 2022-Q3 and
#    shift the previous rows (2022-Q2 / 2022-Q1 / 2021-Q4) down by one.
#    Done via bottom-up Copy (rather than Rows.Insert) so each moved row
#    keeps its original style along with its value.
# ---------------------------------------------------------------------
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A2:D2"))

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 11
$totalSheet.Cells.Item(2, 4).Value = 1.89

# Renumber the index column (A) for the rows that shifted down so it stays
# the contiguous 0,1,2,3 sequence the sheet had before.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
